$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1873536299765808
$ws.Range("C2").Value = 0.5807962529274004
$ws.Range("J2").Value = 0.0117096018735363
$ws.Range("P2").Value = 0.1241217798594848
$ws.Range("S2").Value = 0.09601873536299765
$ws.Range("B3").Value = 0.01568627450980392
$ws.Range("C3").Value = 0.02745098039215686
$ws.Range("J3").Value = 0.04313725490196078
$ws.Range("P3").Value = 0.7333333333333333
$ws.Range("S3").Value = 0.1803921568627451
$ws.Range("J4").Value = 0.02702702702702703
$ws.Range("P4").Value = 0.6486486486486487
$ws.Range("S4").Value = 0.3243243243243243
$ws.Range("B6").Value = 0.06319702602230483
$ws.Range("D6").Value = 0.003717472118959108
$ws.Range("F6").Value = 0.02230483271375465
$ws.Range("J6").Value = 0.3308550185873606
$ws.Range("O6").Value = 0.01858736059479554
$ws.Range("Q6").Value = 0.1226765799256506
$ws.Range("R6").Value = 0.104089219330855
$ws.Range("S6").Value = 0.3345724907063197
$ws.Range("B7").Value = 0.1085972850678733
$ws.Range("D7").Value = 0.02714932126696833
$ws.Range("F7").Value = 0.04977375565610859
$ws.Range("J7").Value = 0.1357466063348416
$ws.Range("O7").Value = 0.01809954751131222
$ws.Range("Q7").Value = 0.2126696832579185
$ws.Range("R7").Value = 0.04524886877828054
$ws.Range("S7").Value = 0.4027149321266968
$ws.Range("B8").Value = 0.1452282157676349
$ws.Range("D8").Value = 0.01037344398340249
$ws.Range("E8").Value = 0.004149377593360996
$ws.Range("F8").Value = 0.07468879668049792
$ws.Range("J8").Value = 0.1078838174273859
$ws.Range("O8").Value = 0.02074688796680498
$ws.Range("Q8").Value = 0.1514522821576763
$ws.Range("R8").Value = 0.07676348547717843
$ws.Range("S8").Value = 0.4087136929460581
$ws.Range("B9").Value = 0.0797872340425532
$ws.Range("D9").Value = 0.01063829787234043
$ws.Range("F9").Value = 0.05319148936170213
$ws.Range("J9").Value = 0.06382978723404255
$ws.Range("O9").Value = 0.02127659574468085
$ws.Range("Q9").Value = 0.2553191489361702
$ws.Range("R9").Value = 0.0797872340425532
$ws.Range("S9").Value = 0.4361702127659575
$ws.Range("B10").Value = 0.136986301369863
$ws.Range("D10").Value = 0.01500326157860404
$ws.Range("F10").Value = 0.07110241356816699
$ws.Range("J10").Value = 0.1304631441617743
$ws.Range("O10").Value = 0.01891715590345727
$ws.Range("Q10").Value = 0.2048271363339856
$ws.Range("R10").Value = 0.08153946510110893
$ws.Range("S10").Value = 0.3411611219830398
$ws.Range("G11").Value = 0.1523341523341523
$ws.Range("J11").Value = 0.14004914004914
$ws.Range("K11").Value = 0.2383292383292383
$ws.Range("L11").Value = 0.457002457002457
$ws.Range("S11").Value = 0.01228501228501228
$ws.Range("G12").Value = 0.6683937823834197
$ws.Range("J12").Value = 0.2538860103626943
$ws.Range("K12").Value = 0.01036269430051814
$ws.Range("L12").Value = 0.0310880829015544
$ws.Range("S12").Value = 0.03626943005181347
$ws.Range("G13").Value = 0.6724137931034483
$ws.Range("J13").Value = 0.2586206896551724
$ws.Range("S13").Value = 0.06896551724137931
$ws.Range("F15").Value = 0.02755905511811024
$ws.Range("H15").Value = 0.1535433070866142
$ws.Range("I15").Value = 0.06299212598425197
$ws.Range("J15").Value = 0.3622047244094488
$ws.Range("K15").Value = 0.08661417322834646
$ws.Range("M15").Value = 0.01574803149606299
$ws.Range("O15").Value = 0.07874015748031496
$ws.Range("S15").Value = 0.2125984251968504
$ws.Range("F16").Value = 0.02362204724409449
$ws.Range("H16").Value = 0.1181102362204724
$ws.Range("I16").Value = 0.07480314960629922
$ws.Range("J16").Value = 0.4015748031496063
$ws.Range("K16").Value = 0.1220472440944882
$ws.Range("M16").Value = 0.01968503937007874
$ws.Range("O16").Value = 0.06692913385826772
$ws.Range("S16").Value = 0.1732283464566929
$ws.Range("F17").Value = 0.017578125
$ws.Range("H17").Value = 0.158203125
$ws.Range("I17").Value = 0.0859375
$ws.Range("J17").Value = 0.42578125
$ws.Range("K17").Value = 0.1328125
$ws.Range("M17").Value = 0.021484375
$ws.Range("N17").Value = 0.001953125
$ws.Range("O17").Value = 0.0546875
$ws.Range("S17").Value = 0.1015625
$ws.Range("F18").Value = 0.04186046511627907
$ws.Range("H18").Value = 0.1906976744186047
$ws.Range("I18").Value = 0.06976744186046512
$ws.Range("J18").Value = 0.413953488372093
$ws.Range("K18").Value = 0.1116279069767442
$ws.Range("M18").Value = 0.009302325581395349
$ws.Range("O18").Value = 0.06976744186046512
$ws.Range("S18").Value = 0.09302325581395349
$ws.Range("F19").Value = 0.02648532569792412
$ws.Range("H19").Value = 0.2083035075161059
$ws.Range("I19").Value = 0.06728704366499642
$ws.Range("J19").Value = 0.3808160343593415
$ws.Range("K19").Value = 0.1145311381531854
$ws.Range("M19").Value = 0.02791696492483894
$ws.Range("O19").Value = 0.06800286327845383
$ws.Range("S19").Value = 0.1066571224051539
